# Applies the "Añadir comentario act viterbi" change: the data rows 2-10
# (POS-tag labels in column A plus their Viterbi probability values in
# columns B:I) are reordered. Row 9 (NCMS000) keeps its place; every other
# row's full content moves to a new row while keeping its own values intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contents for rows 2..10, columns A..I (after the reorder)
$data = @(
    @("NCFS000", 0, 0, 0, 0, 0, 0, 0, 0),
    @("AQ0MS0", 0, 0, 0, 0, 0, 0, 0, 0),
    @("AQ0CS0", 0, 0, 0, 0, 0.006755488834678175, 0, 0, 0),
    @("SPS00", 0, 0.007871568247054057, 0, 0, 0, 0.05138384827938065, 0, 0),
    @("DA0MS0", 0, 0, 0.6456245689480303, 0, 0, 0, 0, 0),
    @("Fp", 0, 0, 0, 0, 0, 0, 0, 0.005586415856583765),
    @("VMIP3S0", 0.001078422914330084, 0, 0, 0, 0, 0, 0, 0),
    @("NCMS000", 0, 0, 0, 0.003542245300487059, 0, 0, 0, 0),
    @("NCMP000", 0, 0, 0, 0, 0, 0, 0.00537109375, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
